$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-27 from 45333 to 45334 (serial date values)
$ws.Range("C2:C27").Value = 45334
